# Update the "Metadata" worksheet of the CodeSystem workbook:
#  - Experimental flag text changes from "false" to "true"
#  - Date value is refreshed to a newer timestamp
#  - Case Sensitive value is now populated with "true"
#
# Note: assigning the literal word true/false straight to Range.Value makes
# Excel auto-convert the cell to a Boolean. The source file stores these as
# plain text, so we build the literal text via a text formula in a scratch
# cell and paste-special the computed value back in (this keeps the cell a
# normal text/string cell instead of a Boolean, and does not disturb the
# existing cell style).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$scratch = $ws.Range("D1")

$scratch.Formula = "=""true"""
$scratch.Copy()
$ws.Range("B7").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("B14").PasteSpecial(-4163) # xlPasteValues
$scratch.ClearContents()

$ws.Range("B8").Value = "2023-02-16T14:43:10-06:00"
